$d = $word.ActiveDocument

# wdReplaceAll = 2
$wdReplaceAll = 2

# --- Simple 1:1 unique-text replacements ---

$d.Content.Find.Execute(
    "Input (SubjectContext): (C, S, P, O)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Input (SubjectContext): (C, SK, PK, OK) / Composite SK(PK, OK) Statement",
    $wdReplaceAll)

$d.Content.Find.Execute(
    "Input (PredicateContext): (C, P, S, O)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Input (PredicateContext): (C, PK, SK, OK) / Composite PK(SK, OK) Mapping",
    $wdReplaceAll)

$d.Content.Find.Execute(
    "Input (ObjectContext): (C, O, P, S)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Input (ObjectContext): (C, OK, PK, SK) / Composite OK(PK, SK): Behavior",
    $wdReplaceAll)

$d.Content.Find.Execute(
    "Output (SubjectContext): (C, S, P, O) / Composite SK(PK, OK) Statement Kinds?", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Output (SubjectContext): (C, SK, PK, OK) / Composite SK(PK, OK) Statement",
    $wdReplaceAll)

$d.Content.Find.Execute(
    "Output (PredicateContext): (C, P, S, O) / Composite PK(SK, OK) Mapping Kinds?", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Output (PredicateContext): (C, PK, SK, OK) / Composite PK(SK, OK) Mapping",
    $wdReplaceAll)

$d.Content.Find.Execute(
    "Output (ObjectContext): (C, O, P, S) / Composite OK(PK, SK) Behavior Kinds?", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Output (ObjectContext): (C, OK, PK, SK) / Composite OK(PK, SK): Behavior",
    $wdReplaceAll)

# --- The two ambiguous "Input: (C, S, P, O)" lines (there are 3 occurrences;
#     only the 2nd and 3rd change, the 1st stays as-is). Target them by
#     paragraph index found via direct text scan. ---

$targets = @(
    @{ Old = "Input: (C, S, P, O)"; New = "Input: (C P, S, O)" },
    @{ Old = "Input: (C, S, P, O)"; New = "Input: (C, O, P, S)" }
)
$seen = 0
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Input: (C, S, P, O)" + [char]13) {
        $seen = $seen + 1
        if ($seen -ge 2) {
            $target = $targets[$seen - 2]
            $rng = $p.Range
            $rng.Find.Execute(
                $target.Old, $true, $false, $false, $false, $false,
                $true, 1, $false, $target.New, $wdReplaceAll)
        }
    }
}

# --- Insert a new empty list paragraph right after the
#     "Output (ObjectContext): ..." paragraph. ---

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Output (ObjectContext): (C, OK, PK, SK) / Composite OK(PK, SK): Behavior" + [char]13) {
        $p.Range.InsertParagraphAfter()
        break
    }
}
